{"js": "// Update each three-digit \u00d7 one-digit multiplication answer cell in the\n// table with its new equation, matched by the cell's current (old) text so\n// that formatting (font/size) on the existing run is preserved.\nconst replacements = {\n  \"652\u00d76=3912\": \"744\u00d76=4464\",\n  \"454\u00d77=3178\": \"712\u00d77=4984\",\n  \"435\u00d73=1305\": \"564\u00d78=4512\",\n  \"487\u00d77=3409\": \"824\u00d72=1648\",\n  \"990\u00d76=5940\": \"975\u00d78=7800\",\n  \"442\u00d74=1768\": \"789\u00d73=2367\",\n  \"579\u00d77=4053\": \"624\u00d76=3744\",\n  \"819\u00d73=2457\": \"535\u00d79=4815\",\n  \"683\u00d72=1366\": \"638\u00d76=3828\",\n  \"892\u00d75=4460\": \"449\u00d79=4041\",\n  \"372\u00d77=2604\": \"316\u00d79=2844\",\n  \"988\u00d76=5928\": \"190\u00d76=1140\",\n  \"581\u00d76=3486\": \"416\u00d72=832\",\n  \"153\u00d75=765\": \"925\u00d78=7400\",\n  \"837\u00d74=3348\": \"617\u00d74=2468\",\n  \"799\u00d74=3196\": \"482\u00d73=1446\",\n  \"502\u00d73=1506\": \"625\u00d76=3750\",\n  \"474\u00d74=1896\": \"462\u00d76=2772\",\n  \"725\u00d73=2175\": \"978\u00d74=3912\",\n  \"198\u00d76=1188\": \"797\u00d76=4782\",\n  \"980\u00d78=7840\": \"169\u00d79=1521\",\n  \"177\u00d76=1062\": \"500\u00d74=2000\",\n  \"994\u00d78=7952\": \"436\u00d76=2616\",\n  \"999\u00d72=1998\": \"206\u00d74=824\",\n  \"605\u00d77=4235\": \"530\u00d76=3180\",\n};\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet appliedCount = 0;\n\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items/text\");\n    await context.sync();\n\n    for (const paragraph of paragraphs.items) {\n      const currentText = paragraph.text;\n      const replacement = replacements[currentText];\n      if (replacement !== undefined) {\n        paragraph.insertText(replacement, \"Replace\");\n        appliedCount++;\n      }\n    }\n  }\n}\n\nawait context.sync();\n\nif (appliedCount !== Object.keys(replacements).length) {\n  throw new Error(\n    `Expected to apply ${Object.keys(replacements).length} replacements, applied ${appliedCount}.`\n  );\n}\n", "ps1": "# Update each three-digit x one-digit multiplication answer cell in the\n# table with its new equation via Find & Replace over the whole document\n# content range (wdFindContinue=1 keeps the search inside Content without\n# prompting to wrap; wdReplaceAll=2 replaces every match of that exact text).\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n  \"652\u00d76=3912\" = \"744\u00d76=4464\"\n  \"454\u00d77=3178\" = \"712\u00d77=4984\"\n  \"435\u00d73=1305\" = \"564\u00d78=4512\"\n  \"487\u00d77=3409\" = \"824\u00d72=1648\"\n  \"990\u00d76=5940\" = \"975\u00d78=7800\"\n  \"442\u00d74=1768\" = \"789\u00d73=2367\"\n  \"579\u00d77=4053\" = \"624\u00d76=3744\"\n  \"819\u00d73=2457\" = \"535\u00d79=4815\"\n  \"683\u00d72=1366\" = \"638\u00d76=3828\"\n  \"892\u00d75=4460\" = \"449\u00d79=4041\"\n  \"372\u00d77=2604\" = \"316\u00d79=2844\"\n  \"988\u00d76=5928\" = \"190\u00d76=1140\"\n  \"581\u00d76=3486\" = \"416\u00d72=832\"\n  \"153\u00d75=765\" = \"925\u00d78=7400\"\n  \"837\u00d74=3348\" = \"617\u00d74=2468\"\n  \"799\u00d74=3196\" = \"482\u00d73=1446\"\n  \"502\u00d73=1506\" = \"625\u00d76=3750\"\n  \"474\u00d74=1896\" = \"462\u00d76=2772\"\n  \"725\u00d73=2175\" = \"978\u00d74=3912\"\n  \"198\u00d76=1188\" = \"797\u00d76=4782\"\n  \"980\u00d78=7840\" = \"169\u00d79=1521\"\n  \"177\u00d76=1062\" = \"500\u00d74=2000\"\n  \"994\u00d78=7952\" = \"436\u00d76=2616\"\n  \"999\u00d72=1998\" = \"206\u00d74=824\"\n  \"605\u00d77=4235\" = \"530\u00d76=3180\"\n}\n\n$totalReplacements = 0\nforeach ($oldText in $replacements.Keys) {\n  $newText = $replacements[$oldText]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Find.Execute could not locate expected text: $oldText\"\n  }\n  $totalReplacements++\n}\n\nWrite-Output \"Applied $totalReplacements replacements.\"\n"}
